$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new data row (row 2) with attendeeId, campId and isWithdrawn values,
# matching the existing header row columns: attendeeId | campId | isWithdrawn
$ws.Range("A2").Value = "C133313"
$ws.Range("B2").Value = "abcde"
$ws.Range("C2").Value = $false
